$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 170, shifting existing rows 170-268 down to 172-270.
$ws.Rows.Item(170).Resize(2).Insert()

# New row 170: new weekly record (date 2022-07-25 / serial 44767)
$ws.Range("A170").Value = 1
$ws.Range("B170").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C170").Value = "Arica y Parinacota"
$ws.Range("D170").Value = 44767
$ws.Range("E170").Value = 15
$ws.Range("F170").Value = "Fruta"
$ws.Range("G170").Value = 100108
$ws.Range("H170").Value = "Tropicales y subtropicales"
$ws.Range("I170").Value = 100108006
$ws.Range("J170").Value = "Plátano"
$ws.Range("K170").Value = "Sin especificar"
$ws.Range("L170").Value = "Pintón"
$ws.Range("M170").Value = 120
$ws.Range("N170").Value = 14000
$ws.Range("O170").Value = 15000
$ws.Range("P170").Value = 14500
$ws.Range("Q170").Value = "$/caja 20 kilos"
$ws.Range("R170").Value = "Bolivia"
$ws.Range("S170").Value = 725
$ws.Range("T170").Value = 20

# New row 171: new weekly record (date 2022-07-25 / serial 44767)
$ws.Range("A171").Value = 1
$ws.Range("B171").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C171").Value = "Arica y Parinacota"
$ws.Range("D171").Value = 44767
$ws.Range("E171").Value = 15
$ws.Range("F171").Value = "Fruta"
$ws.Range("G171").Value = 100108
$ws.Range("H171").Value = "Tropicales y subtropicales"
$ws.Range("I171").Value = 100108006
$ws.Range("J171").Value = "Plátano"
$ws.Range("K171").Value = "Sin especificar"
$ws.Range("L171").Value = "Pintón"
$ws.Range("M171").Value = 120
$ws.Range("N171").Value = 28000
$ws.Range("O171").Value = 29000
$ws.Range("P171").Value = 28500
$ws.Range("Q171").Value = "$/caja 20 kilos"
$ws.Range("R171").Value = "Ecuador"
$ws.Range("S171").Value = 1425
$ws.Range("T171").Value = 20
